$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.391.39"
$ws.Range("E2").Value = "  -5.00%  "
$ws.Range("D3").Value = "3.250.72"
$ws.Range("E3").Value = "  -8.27%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'586.74"
$ws.Range("E5").Value = "  -5.06%  "
$ws.Range("D6").Value = "'153.37"
$ws.Range("E6").Value = "  -12.38%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "3.241.65"
$ws.Range("E8").Value = "  -8.36%  "
$ws.Range("D9").Value = "'0.545"
$ws.Range("E9").Value = "  -11.07%  "
$ws.Range("D10").Value = "'0.173"
$ws.Range("E10").Value = "  -13.40%  "
$ws.Range("D11").Value = "'6.83"
$ws.Range("E11").Value = "  -5.66%  "
$ws.Range("D12").Value = "'0.509"
$ws.Range("E12").Value = "  -13.60%  "
$ws.Range("D13").Value = "'38.75"
$ws.Range("E13").Value = "  -17.17%  "
$ws.Range("D14").Value = "'0.0000245"
$ws.Range("E14").Value = "  -11.50%  "
$ws.Range("D15").Value = "3.772.20"
$ws.Range("E15").Value = "  -8.33%  "
$ws.Range("D16").Value = "67.514.52"
$ws.Range("E16").Value = "  -4.94%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.257.50"
$ws.Range("E17").Value = "  -8.16%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'545.25"
$ws.Range("E18").Value = "  -10.97%  "
$ws.Range("D19").Value = "'7.27"
$ws.Range("E19").Value = "  -14.01%  "
$ws.Range("E20").Value = "  -5.87%  "
$ws.Range("D21").Value = "'15.29"
$ws.Range("E21").Value = "  -14.11%  "
$ws.Range("D22").Value = "'0.768"
$ws.Range("E22").Value = "  -13.63%  "
$ws.Range("D23").Value = "'7.85"
$ws.Range("E23").Value = "  -13.59%  "
$ws.Range("D24").Value = "'85.80"
$ws.Range("E24").Value = "  -12.99%  "
$ws.Range("D25").Value = "'13.55"
$ws.Range("E25").Value = "  -13.96%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'3.20"
$ws.Range("E27").Value = "  -15.83%  "
$ws.Range("D28").Value = "'8.24"
$ws.Range("E28").Value = "  -10.09%  "
$ws.Range("D29").Value = "'29.58"
$ws.Range("E29").Value = "  -12.79%  "
$ws.Range("E30").Value = "  -17.80%  "
$ws.Range("D31").Value = "'2.73"
$ws.Range("E31").Value = "  -10.53%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("E32").Value = "  -11.42%  "
$ws.Range("D33").Value = "'553.11"
$ws.Range("E33").Value = "  -12.54%  "
$ws.Range("D34").Value = "'6.62"
$ws.Range("E34").Value = "  -19.10%  "
$ws.Range("D35").Value = "'5.79"
$ws.Range("E35").Value = "  -15.93%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'53.89"
$ws.Range("E37").Value = "  -5.50%  "
$ws.Range("D38").Value = "'0.0439"
$ws.Range("E38").Value = "  -8.28%  "
$ws.Range("D39").Value = "'9.26"
$ws.Range("E39").Value = "  -14.80%  "
$ws.Range("D40").Value = "'0.0852"
$ws.Range("E40").Value = "  -15.04%  "
$ws.Range("D41").Value = "'0.128"
$ws.Range("E41").Value = "  -11.60%  "
$ws.Range("D42").Value = "2.946.39"
$ws.Range("E42").Value = "  -12.67%  "
$ws.Range("D43").Value = "'2.62"
$ws.Range("E43").Value = "  -25.32%  "
$ws.Range("D44").Value = "'0.263"
$ws.Range("E44").Value = "  -16.28%  "
$ws.Range("D45").Value = "0.0₃0588"
$ws.Range("E45").Value = "  -20.96%  "
$ws.Range("D46").Value = "'2.38"
$ws.Range("E46").Value = "  -20.64%  "
$ws.Range("D47").Value = "'26.42"
$ws.Range("E47").Value = "  -18.07%  "
$ws.Range("D48").Value = "'2.14"
$ws.Range("E48").Value = "  -16.90%  "
$ws.Range("D50").Value = "'125.55"
$ws.Range("E50").Value = "  -5.72%  "
$ws.Range("E51").Value = "  -12.81%  "
